$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.69%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.78%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.104"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.64%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07825"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.59%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.315"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-9.83%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.34%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.837"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.21%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9148"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.83%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1755"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.29%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07593"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.81%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09123"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.86%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03090"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.42%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.56%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'1.26%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005919"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.02%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.489"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.248"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.28%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3290"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.21%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.85%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.035"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-12.40%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1791"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'12.02%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04592"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.31%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001251"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.59%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004458"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.03%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'6.07%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-1.30%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01772"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.43%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.20%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007252"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.87%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1359"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.14%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002192"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.10%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-1.55%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006191"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.85%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.16%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'28.94%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7453"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-9.17%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.16%"
$ws.Range("E50").Style = "Normal"
